# Update the "drift index" values on Sheet1 with unrounded, higher-precision
# figures (prep work before removing the Stephens-MacCall filtering column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    "B2" = 0.44260906702452502
    "C2" = 0.44931330993184498
    "D2" = 0.76410993471033894
    "E2" = 0.84954484595345803

    "B3" = 0.134386609473158
    "C3" = 0.14154159191512899
    "D3" = 0.16107345657981301
    "E3" = 0.33249135133184399

    "B4" = 0.24156858465198
    "C4" = 0.23992989606630999
    "D4" = 0.86528804024849104
    "E4" = 0.91618807603727404

    "B5" = 0.31966533038908101
    "C5" = 0.30116402820974297
    "D5" = 0.448118735029992
    "E5" = 0.20871137050566399

    "B6" = 0.17854207131806199
    "C6" = 0.183113197028901
    "D6" = 0.25622049823083798
    "E6" = 0.25351899448630599

    "B7" = 0.15191200095070601
    "C7" = 0.178188399492129
    "D7" = 0.42244511602686402
    "E7" = 0.50878893020593696
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}

# Match the recorded selection from the saved file: B2:E7 active at B2.
$ws.Range("B2:E7").Select()
